# Update the "Förändrad" (Changed) date column (C) from 2023-09-09 (45178)
# to 2023-09-10 (45179) for every data row (rows 2 through 398).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 398
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = 45179
}
